$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("function")

# Insert two new rows above row 10, shifting existing rows (10-22) down to (12-24)
$ws.Rows("10:11").Insert() | Out-Null

# The "dataType" / "41200~41300" row that used to be row 17 is now at row 19.
# Move its content up into the newly inserted row 10.
$a19 = $ws.Range("A19").Value2
$c19 = $ws.Range("C19").Value2
$ws.Range("A10").Value = $a19
$ws.Range("C10").Value = $c19

# Clear the now-duplicated old location.
$ws.Range("A19").ClearContents() | Out-Null
$ws.Range("C19").ClearContents() | Out-Null

# Fill the second newly inserted row (11) with the new "sms" entry.
$ws.Range("A11").Value = "sms"
$ws.Range("C11").Value = "41300~41400"

# Make "function" sheet the active/selected sheet with J13 as the active cell
# (this clears tabSelected on whichever sheet had it before, e.g. "db").
$ws.Select() | Out-Null
$ws.Range("J13").Select() | Out-Null
